$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "98×96=9408" "40×31=1240"
Replace-Text "93×71=6603" "75×49=3675"
Replace-Text "62×43=2666" "22×84=1848"
Replace-Text "43×44=1892" "48×19=912"
Replace-Text "42×59=2478" "82×75=6150"
Replace-Text "62×67=4154" "70×75=5250"
Replace-Text "46×11=506" "37×98=3626"
Replace-Text "42×51=2142" "34×91=3094"
Replace-Text "77×98=7546" "67×15=1005"
Replace-Text "61×81=4941" "59×62=3658"
Replace-Text "72×91=6552" "94×79=7426"
Replace-Text "33×25=825" "80×48=3840"
Replace-Text "30×87=2610" "66×23=1518"
Replace-Text "72×88=6336" "17×52=884"
Replace-Text "26×17=442" "94×14=1316"
Replace-Text "97×95=9215" "89×41=3649"
Replace-Text "27×37=999" "44×66=2904"
Replace-Text "75×60=4500" "95×65=6175"
Replace-Text "18×39=702" "37×68=2516"
Replace-Text "29×75=2175" "88×82=7216"
Replace-Text "60×83=4980" "49×47=2303"
Replace-Text "69×93=6417" "40×66=2640"
Replace-Text "76×77=5852" "30×36=1080"
Replace-Text "59×83=4897" "14×81=1134"
Replace-Text "97×37=3589" "86×32=2752"

Write-Output "Done replacing"
